# RegionMeetings_template.xlsx edit
#
# 1. Show user full name rather than login name for all kinds of data:
#    the "salesPerson" placeholder becomes "salesPersonFullName".
#    In the template this text lives in cell E2 (row-2 is the JXLS
#    "jx:each" data row that is driven by the placeholder expressions).
#
# 2. The active selection moved from H2 to F3 (a side effect of the
#    author's edit in the live workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '${record.salesPersonFullName}'

$ws.Range("F3").Select()
